# Apply the "RegisterPageTest" data-sheet edit:
#  - rename Sheet3 -> Register and populate it with the register-page test data
#  - add a missing data column (D2) on the Search sheet
#  - move the active tab / selection to match the new authoring state

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Search
$ws2 = $wb.Worksheets.Item(2)   # Login
$ws3 = $wb.Worksheets.Item(3)   # Sheet3 -> Register

# ---------------------------------------------------------------------------
# Sheet3 -> "Register"
# ---------------------------------------------------------------------------
$ws3.Name = "Register"

# Column widths (matches the authored <cols> block). The host engine snaps
# ColumnWidth to 1/6-character increments, so the inputs below are chosen to
# land on the closest reproducible width to the authored value.
$ws3.Columns.Item(2).ColumnWidth  = 17.15  # -> 18
$ws3.Columns.Item(3).ColumnWidth  = 15.0   # -> 15.8333 (target 15.85546875)
$ws3.Columns.Item(4).ColumnWidth  = 12.65  # -> 13.5    (target 13.5703125)
$ws3.Columns.Item(5).ColumnWidth  = 20.65  # -> 21.5    (target 21.42578125)
$ws3.Columns.Item(6).ColumnWidth  = 18.65  # -> 19.5    (target 19.42578125)
$ws3.Columns.Item(7).ColumnWidth  = 16.85  # -> 17.6667 (target 17.7109375)
$ws3.Columns.Item(8).ColumnWidth  = 12.65  # -> 13.5    (target 13.5703125)
$ws3.Columns.Item(9).ColumnWidth  = 13.65  # -> 14.5    (target 14.5703125)
$ws3.Columns.Item(10).ColumnWidth = 15.3   # -> 16.1667 (target 16.140625)
$ws3.Columns.Item(11).ColumnWidth = 10.5   # -> 11.3333 (target 11.28515625)

# Header row (row 1) -- written in the exact order the strings were first
# authored so the shared-string table comes out in the original order.
$ws3.Range("A1").Value = "TC_ID"
$ws3.Range("B1").Value = "password_invalid"
$ws3.Range("C1").Value = "password_valid"
$ws3.Range("D1").Value = "fullName"
$ws3.Range("F1").Value = "mobileNumber_valid"
$ws3.Range("E1").Value = "mobileNumber_invalid"
$ws3.Range("G1").Value = "designation"
$ws3.Range("H1").Value = "basicEducation"
$ws3.Range("I1").Value = "filePath_invalid"
$ws3.Range("J1").Value = "filePath_valid"
$ws3.Range("K1").Value = "resumeText"
$ws3.Range("L1").Value = "specialization"
$ws3.Range("M1").Value = "institution"
$ws3.Range("N1").Value = "yearOfCompletion"
$ws3.Range("O1").Value = "currentIndustry"
$ws3.Range("P1").Value = "functionalArea"
$ws3.Range("Q1").Value = "role"
$ws3.Range("R1").Value = "salary"
$ws3.Range("S1").Value = "jobAlertName"

# Data row (row 2) -- again in first-authored order so the shared string
# table lands with "mukeshkumar" (D2) last, matching the source commit.
$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = "aks"
$ws3.Range("C2").Value = "godisking"
$ws3.Range("E2").Value = 12345
$ws3.Range("F2").Value = 9900224430
$ws3.Range("G2").Value = "Senior Software Engineer"
$ws3.Range("H2").Value = "B.Tech/B.E."
$ws3.Range("I2").Value = "C:\\Users\\Admin\\Desktop\\resume.txt"
$ws3.Range("J2").Value = "C:\\Users\\Admin\\Desktop\\resume.docx"
$ws3.Range("K2").Value = "My Resume"
$ws3.Range("L2").Value = "Computers"
$ws3.Range("M2").Value = "Anna University"
$ws3.Range("N2").Value = 2005
$ws3.Range("O2").Value = "IT-Software/Software Services"
$ws3.Range("P2").Value = "IT Software - Application Programming / Maintenance"
$ws3.Range("Q2").Value = "Testing Engnr"
$ws3.Range("R2").Value = 6
$ws3.Range("S2").Value = "MyAlert"
$ws3.Range("D2").Value = "mukeshkumar"

# ---------------------------------------------------------------------------
# Sheet1 ("Search") -- new D2 value + selection change
# ---------------------------------------------------------------------------
$ws1.Range("D2").Value = 6
$ws1.Range("D5").Select()

# ---------------------------------------------------------------------------
# Activate Register last so it becomes the active tab / selected sheet.
# ---------------------------------------------------------------------------
$ws3.Range("F4").Select()
$ws3.Activate()
